$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "325.05"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.24%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.64"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.88%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.654"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "6.57%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08036"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.81%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.018"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.80%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.487"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.71%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.618"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.40%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.933"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.47%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9223"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.09%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1243"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-8.63%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1964"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.40%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.731"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "21.42%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09186"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.19%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03560"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.09%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.1049"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "9.42%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001295"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.12%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006119"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.82%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.41%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.91%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1371"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.88%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04378"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.09%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001263"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.33%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004606"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.69%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "2.54%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02521"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.66%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05330"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.47%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007460"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.47%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009914"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.63%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1404"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.57%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.03%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01108"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "11.72%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006689"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.53%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.02%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002281"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-5.06%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.02%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.02%"
